# Milestone guidelines.docx — record the completion date/time for the
# "hospital sounds" milestone row (Day-2 table row whose notes cell
# mentions "Find and implement Hospital audio.").
#
# Before: that row's "Date Completed: " paragraph is empty (just the
#         label), sitting right before the document's _GoBack bookmark.
# After:  the label paragraph gets "16/08/10" appended as its own bold
#         run (and its paragraph mark picks up bold too), and a brand
#         new centered/bold "Time: 4:35 Pm" paragraph is inserted right
#         after it, carrying the _GoBack bookmark forward onto itself.

$d = $word.ActiveDocument
$paragraphs = $d.Paragraphs

# Anchor on the milestone row that talks about hospital audio, per the
# commit message ("added hospital sounds").
$hospitalIndex = -1
for ($i = 1; $i -le $paragraphs.Count; $i++) {
    if ($paragraphs.Item($i).Range.Text -like "*Hospital*") {
        $hospitalIndex = $i
    }
}

# From there, walk forward to the next still-blank "Date Completed: "
# paragraph — that's the one belonging to this same table row.
$targetIndex = -1
for ($i = $hospitalIndex; $i -le $paragraphs.Count; $i++) {
    $t = $paragraphs.Item($i).Range.Text
    # Range.Text carries a trailing paragraph-mark (and, since this is
    # the last paragraph in its table cell, a cell-mark) control char.
    $trimmed = $t.TrimEnd([char]13, [char]7)
    if ($trimmed -eq "Date Completed: ") {
        $targetIndex = $i
        break
    }
}

$p = $paragraphs.Item($targetIndex)
$r = $p.Range

# Replace this single paragraph with the two finished paragraphs
# (keeping the run split / bold formatting / bookmark placement exactly
# as Word would produce them), in one shot so no stray empty paragraph
# or merge step is needed.
$xml = '<w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:sz w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve">Date Completed: </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="28"/></w:rPr><w:t>16/08/10</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:sz w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve">Time: </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="28"/></w:rPr><w:t>4:35 Pm</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$r.InsertXML($xml)
